$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A values
$ws.Range("A1").Value = "allumer l'univers 2 mais reste allumer"
$ws.Range("A2").Value = "allumer l'univers 2 mais reste allumer"
$ws.Range("A3").Value = "allumer l'universa 2 et 4 reste allumer"
$ws.Range("A4").Value = "allumer unvivers 2 et 4 reste allumer"
$ws.Range("A5").Value = "allumer unvers 2 reste allumer"
$ws.Range("A6").Value = "allumer univers 2 reste allumer"
$ws.Range("A7").Value = "meme"
$ws.Range("A8").Value = "meme"
$ws.Range("A11").Value = "Tout les tests marchent sauf que sur case chemine les leds reste allumer mais pas sur toit-vitre"
$ws.Range("A12").Value = "pour les fad in fad out et parabole. l'animation est bizard elle reste beaucoup de temps allumer(pas sur de pouvoir employer tout les valeurs de 0 à 255)"
$ws.Range("A13").Value = "pour la luminosité à partir de 15 16 on ne voit plus d'augmentation de la luminosité ou tres tres legerement je n'arrive pas à le voir à l'œil nu"

# Column E values
$ws.Range("E1").Value = "test1 "
$ws.Range("E2").Value = "test1 sync"
$ws.Range("E3").Value = "test2"
$ws.Range("E4").Value = "test 2 sync"
$ws.Range("E5").Value = "test1 continu"
$ws.Range("E6").Value = "test1 continu sync"
$ws.Range("E7").Value = "test 2 continue "
$ws.Range("E8").Value = "test 2 continue  sync"

# Column F value
$ws.Range("F8").Value = " "

# Column I value
$ws.Range("I4").Value = "pas de différence entre le art sync et le non artsync"

# Selection in the saved file points at I4
$ws.Range("I4").Select()
